$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_8a_Links")

# Row 24: update German link URL, clear English link
$ws.Range("E24").Value = "https://www.destatis.de/DE/Themen/Arbeit/Verdienste/Verdienste-GenderPayGap/Tabellen/ugpg-02-bundeslaender-ab-2014.html"
$ws.Range("F24").Value = ""

# Row 37: update German link URL, clear English link
$ws.Range("E37").Value = "https://www.destatis.de/DE/Themen/Arbeit/Verdienste/Verdienste-GenderPayGap/Tabellen/ugpg-01-gebietsstand.html"
$ws.Range("F37").Value = ""

# Row 47: update source code
$ws.Range("B47").Value = "Q_EUROPEANCOMMISSION"

# Row 77: update source code
$ws.Range("B77").Value = "Q_UBA_1"
